$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Orchid  cosmetics" cell: split the two runs into three runs and wrap
#    them with spell/gram proofErr markers (Orchid / "  " / cosmetics).
# ---------------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("Orchid  cosmetics", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find 'Orchid  cosmetics' text"
}

$p1Xml = '<w:p w14:paraId="1E4257C9" w14:textId="236FFB00" w:rsidR="008C429C" w:rsidRPr="00796EAD" w:rsidRDefault="005D459F" w:rsidP="008C429C"><w:pPr><w:spacing w:before="120" w:after="120"/><w:rPr><w:sz w:val="20"/><w:lang w:val="es-PA"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:eastAsia="SimSun" w:cstheme="minorHAnsi"/><w:color w:val="0070C0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="es-ES" w:eastAsia="zh-CN"/></w:rPr><w:t>Orchid</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:eastAsia="SimSun" w:cstheme="minorHAnsi"/><w:color w:val="0070C0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="es-ES" w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve">  </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="009D0AD1"><w:rPr><w:rFonts w:eastAsia="SimSun" w:cstheme="minorHAnsi"/><w:color w:val="0070C0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="es-ES" w:eastAsia="zh-CN"/></w:rPr><w:t>cosmetics</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/></w:p>'

$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' + $p1Xml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# 2) "Queda excluido el soporte t\u00e9cnico post-proyecto." paragraph: split
#    the single run into three runs, wrapping "post-proyecto" with
#    spellStart/spellEnd proofErr markers.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute("Queda excluido el soporte técnico post-proyecto.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find 'Queda excluido el soporte tecnico post-proyecto.' text"
}

$p2Xml = '<w:p w14:paraId="7B1166C3" w14:textId="77777777" w:rsidR="005D30CE" w:rsidRPr="00972711" w:rsidRDefault="005D30CE" w:rsidP="005D30CE"><w:pPr><w:pStyle w:val="infoblue"/><w:ind w:left="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:i w:val="0"/><w:iCs w:val="0"/><w:color w:val="0070C0"/><w:sz w:val="20"/><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r w:rsidRPr="00972711"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:i w:val="0"/><w:iCs w:val="0"/><w:color w:val="0070C0"/><w:sz w:val="20"/><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve">Queda excluido el soporte técnico </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:i w:val="0"/><w:iCs w:val="0"/><w:color w:val="0070C0"/><w:sz w:val="20"/><w:lang w:val="es-ES"/></w:rPr><w:t>post-proyecto</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:i w:val="0"/><w:iCs w:val="0"/><w:color w:val="0070C0"/><w:sz w:val="20"/><w:lang w:val="es-ES"/></w:rPr><w:t>.</w:t></w:r></w:p>'

$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' + $p2Xml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r2.InsertXML($xml2)

# ---------------------------------------------------------------------------
# 3) Footer: replace the "NUMPAGES" w:fldSimple field with an explicit
#    fldChar begin/instrText/separate/result/end run sequence.
# ---------------------------------------------------------------------------
$ftr = $d.Sections(1).Footers(1)
$p3 = $ftr.Range.Paragraphs(1)
$r3 = $p3.Range

$p3Xml = '<w:p w14:paraId="2C9BA461" w14:textId="5437F5CF" w:rsidR="00837F2F" w:rsidRDefault="00837F2F" w:rsidP="00837F2F"><w:pPr><w:pStyle w:val="Piedepgina"/><w:pBdr><w:top w:val="single" w:sz="4" w:space="1" w:color="auto"/></w:pBdr><w:tabs><w:tab w:val="clear" w:pos="4419"/><w:tab w:val="clear" w:pos="8838"/><w:tab w:val="center" w:pos="5400"/><w:tab w:val="right" w:pos="10800"/></w:tabs><w:spacing w:afterAutospacing="0"/></w:pPr><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r w:rsidR="00A45DF9"><w:tab/></w:r><w:r w:rsidR="00A45DF9"><w:tab/></w:r><w:r w:rsidR="00A45DF9"><w:tab/></w:r><w:r><w:t xml:space="preserve">Página </w:t></w:r><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> PAGE   \* MERGEFORMAT </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="00037805"><w:rPr><w:noProof/></w:rPr><w:t>1</w:t></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:fldChar w:fldCharType="end"/></w:r><w:r w:rsidR="00A45DF9"><w:t xml:space="preserve"> de </w:t></w:r><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> NUMPAGES   \* MERGEFORMAT </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:t>1</w:t></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p>'

$xml3 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/footer1.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.footer+xml"><pkg:xmlData><w:ftr xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + $p3Xml + '</w:ftr></pkg:xmlData></pkg:part></pkg:package>'

$r3.InsertXML($xml3)

Write-Host "Done."
